# Auto-generated Excel COM-interop edit script
# Applies numeric updates to columns H-N across multiple rows/sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1873.5714
$ws.Range("I88").Value = 722
$ws.Range("J88").Value = 2449.3572
$ws.Range("K88").Value = 722
$ws.Range("L88").Value = 2449.3572
$ws.Range("M88").Value = -316
$ws.Range("N88").Value = -3261.3572

$ws.Range("H91").Value = 1873.5714
$ws.Range("I91").Value = 722
$ws.Range("J91").Value = 2449.3572
$ws.Range("K91").Value = 722
$ws.Range("L91").Value = 2449.3572
$ws.Range("M91").Value = 682
$ws.Range("N91").Value = -5257.3572

$ws.Range("H92").Value = 1060.65
$ws.Range("I92").Value = 277
$ws.Range("J92").Value = 5501.3335
$ws.Range("K92").Value = 277
$ws.Range("L92").Value = 5501.3335
$ws.Range("M92").Value = 971
$ws.Range("N92").Value = -7997.3335

$ws.Range("H94").Value = 14662.5
$ws.Range("I94").Value = 3628.5715
$ws.Range("J94").Value = 23244.445
$ws.Range("K94").Value = 3628.5715
$ws.Range("L94").Value = 23244.445
$ws.Range("M94").Value = -3177.5715
$ws.Range("N94").Value = -24146.445

$ws.Range("H96").Value = 269
$ws.Range("I96").Value = 221.5
$ws.Range("K96").Value = 664.5
$ws.Range("M96").Value = 708.5

$ws.Range("H100").Value = 5873.25
$ws.Range("I100").Value = 2396
$ws.Range("J100").Value = 11668.667
$ws.Range("K100").Value = 2396
$ws.Range("L100").Value = 11668.667
$ws.Range("M100").Value = -1855
$ws.Range("N100").Value = -12750.667

$ws.Range("H116").Value = 3555.3447
$ws.Range("I116").Value = 2130.25
$ws.Range("J116").Value = 6722.222
$ws.Range("K116").Value = 2130.25
$ws.Range("L116").Value = 6722.222
$ws.Range("M116").Value = 1311.75
$ws.Range("N116").Value = -13606.222

$ws.Range("H125").Value = 1792.7142
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H132").Value = 3415.3333
$ws.Range("I132").Value = 2251.9167
$ws.Range("J132").Value = 5742.1665
$ws.Range("K132").Value = 6755.750100000001
$ws.Range("L132").Value = 17226.4995
$ws.Range("M132").Value = -4225.750100000001
$ws.Range("N132").Value = -22286.4995

$ws.Range("H138").Value = 23394006
$ws.Range("I138").Value = 83335410
$ws.Range("J138").Value = 7409632.5
$ws.Range("K138").Value = 250006230
$ws.Range("L138").Value = 22228897.5
$ws.Range("M138").Value = -250001090
$ws.Range("N138").Value = -22239177.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20483.285
$ws.Range("I32").Value = 4412.1963
$ws.Range("J32").Value = 149052
$ws.Range("K32").Value = 4412.1963
$ws.Range("L32").Value = 149052
$ws.Range("M32").Value = -4125.1963
$ws.Range("N32").Value = -149626

$ws.Range("H74").Value = 963.0303
$ws.Range("I74").Value = 1027.037
$ws.Range("J74").Value = 675
$ws.Range("K74").Value = 1027.037
$ws.Range("L74").Value = 675
$ws.Range("M74").Value = -153.037
$ws.Range("N74").Value = -2423

$ws.Range("H77").Value = 963.0303
$ws.Range("I77").Value = 1027.037
$ws.Range("J77").Value = 675
$ws.Range("K77").Value = 5135.185
$ws.Range("L77").Value = 3375
$ws.Range("M77").Value = -767.1850000000004
$ws.Range("N77").Value = -12111

$ws.Range("H97").Value = 3910.476
$ws.Range("J97").Value = 6700
$ws.Range("L97").Value = 6700
$ws.Range("N97").Value = -7692

$ws.Range("H102").Value = 111114010
$ws.Range("I102").Value = 3051.6667
$ws.Range("J102").Value = 333335940
$ws.Range("K102").Value = 3051.6667
$ws.Range("L102").Value = 333335940
$ws.Range("M102").Value = -1429.6667
$ws.Range("N102").Value = -333339184

$ws.Range("H110").Value = 4839.4
$ws.Range("I110").Value = 4430.077
$ws.Range("J110").Value = 7500
$ws.Range("K110").Value = 4430.077
$ws.Range("L110").Value = 7500
$ws.Range("M110").Value = -2385.077
$ws.Range("N110").Value = -11590

$ws.Range("H122").Value = 1209.5588
$ws.Range("I122").Value = 1015.4091
$ws.Range("J122").Value = 1565.5
$ws.Range("K122").Value = 3046.2273
$ws.Range("L122").Value = 4696.5
$ws.Range("M122").Value = -596.2273
$ws.Range("N122").Value = -9596.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 956.9048
$ws.Range("I94").Value = 863.9286
$ws.Range("J94").Value = 1142.8572
$ws.Range("K94").Value = 863.9286
$ws.Range("L94").Value = 1142.8572
$ws.Range("M94").Value = -412.9286
$ws.Range("N94").Value = -2044.8572

$ws.Range("H107").Value = 7238.885
$ws.Range("I107").Value = 7509.591
$ws.Range("J107").Value = 5750
$ws.Range("K107").Value = 7509.591
$ws.Range("L107").Value = 5750
$ws.Range("M107").Value = -5589.591
$ws.Range("N107").Value = -9590

$ws.Range("H134").Value = 81153.55499999999
$ws.Range("I134").Value = 105590.86
$ws.Range("J134").Value = 2411.111
$ws.Range("K134").Value = 316772.58
$ws.Range("L134").Value = 7233.333
$ws.Range("M134").Value = -314237.58
$ws.Range("N134").Value = -12303.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 50475
$ws.Range("J53").Value = 50475
$ws.Range("L53").Value = 50475
$ws.Range("N53").Value = -51689

$ws.Range("H99").Value = 1379
$ws.Range("I99").Value = 1367.6666
$ws.Range("J99").Value = 1401.6666
$ws.Range("K99").Value = 1367.6666
$ws.Range("L99").Value = 1401.6666
$ws.Range("M99").Value = 130.3334
$ws.Range("N99").Value = -4397.6666

$ws.Range("H111").Value = 29834
$ws.Range("J111").Value = 29834
$ws.Range("L111").Value = 29834
$ws.Range("N111").Value = -38014

$ws.Range("H118").Value = 24666.666
$ws.Range("J118").Value = 24666.666
$ws.Range("L118").Value = 24666.666
$ws.Range("N118").Value = -27980.666

$ws.Range("H126").Value = 1379
$ws.Range("I126").Value = 1367.6666
$ws.Range("J126").Value = 1401.6666
$ws.Range("K126").Value = 4102.9998
$ws.Range("L126").Value = 4204.9998
$ws.Range("M126").Value = -1632.9998
$ws.Range("N126").Value = -9144.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 653.7368
$ws.Range("I113").Value = 675.8
$ws.Range("J113").Value = 645.8570999999999
$ws.Range("K113").Value = 2027.4
$ws.Range("L113").Value = 1937.5713
$ws.Range("M113").Value = 142.6000000000001
$ws.Range("N113").Value = -6277.5713

$ws.Range("H132").Value = 732039.9399999999
$ws.Range("I132").Value = 1012989.94
$ws.Range("J132").Value = 1570
$ws.Range("K132").Value = 9116909.459999999
$ws.Range("L132").Value = 14130
$ws.Range("M132").Value = -9114379.459999999
$ws.Range("N132").Value = -19190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 44374.332
$ws.Range("J127").Value = 44374.332
$ws.Range("L127").Value = 44374.332
$ws.Range("N127").Value = -54294.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 866.7059
$ws.Range("I22").Value = 678.3077
$ws.Range("J22").Value = 983.3333
$ws.Range("K22").Value = 678.3077
$ws.Range("L22").Value = 983.3333
$ws.Range("M22").Value = -383.3077
$ws.Range("N22").Value = -1573.3333

$ws.Range("H27").Value = 866.7059
$ws.Range("I27").Value = 678.3077
$ws.Range("J27").Value = 983.3333
$ws.Range("K27").Value = 678.3077
$ws.Range("L27").Value = 983.3333
$ws.Range("M27").Value = -571.3077
$ws.Range("N27").Value = -1197.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1111721
$ws.Range("I100").Value = 664.3333
$ws.Range("J100").Value = 3333834.2
$ws.Range("K100").Value = 1328.6666
$ws.Range("L100").Value = 6667668.4
$ws.Range("M100").Value = -787.6666
$ws.Range("N100").Value = -6668750.4

$ws.Range("H107").Value = 1936.1333
$ws.Range("I107").Value = 1505.25
$ws.Range("J107").Value = 2428.5715
$ws.Range("K107").Value = 4515.75
$ws.Range("L107").Value = 7285.7145
$ws.Range("M107").Value = -2595.75
$ws.Range("N107").Value = -11125.7145

$ws.Range("H127").Value = 53618.375
$ws.Range("J127").Value = 53618.375
$ws.Range("L127").Value = 53618.375
$ws.Range("N127").Value = -63538.375
